{"js": "// Office.js (Word JavaScript API) script.\n// Applies the content changes described by the diff:\n//  1. Typo fix: \"Percieved\" -> \"Perceived\" (Skala opazanog stresa / Perceived Stress Scale)\n//  2. Remove \" (Pravilo za 7)\" from the \"P6 (Pravilo za 7) - \" rule label.\n//  3. In rule P7:\n//       \"posljednjih 6 mjeseci\" -> \"posljednja 3 mjeseca\"\n//       \"osim niskog do umjerenog stresa, kao ni predispozicije za oboljenjima\"\n//           -> \"osim stresa i/ili predispozicije ka oboljenjima\"\n//       append \" i korisnik ce se testirati za druge poremecaje iako stres nije visok\"\n//       right after the C3 sentence \"...duzi period)\".\n\nasync function replaceOnce(context, searchText, replacement, insertLocation) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, insertLocation || Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. Fix the \"Percieved\" -> \"Perceived\" typo.\nawait replaceOnce(context, \"Percieved\", \"Perceived\");\n\n// 2. Drop the \" (Pravilo za 7)\" aside after \"P6\".\nawait replaceOnce(context, \" (Pravilo za 7)\", \"\");\n\n// 3a. Tighten the recency window in rule P7.\nawait replaceOnce(context, \"posljednjih 6 mjeseci\", \"posljednja 3 mjeseca\");\n\n// 3b. Reword the exclusion clause in rule P7.\nawait replaceOnce(\n  context,\n  \"osim niskog do umjerenog stresa, kao ni predispozicije za oboljenjima\",\n  \"osim stresa i/ili predispozicije ka oboljenjima\"\n);\n\n// 3c. Append the new trailing clause onto the end of the C3 sentence.\nawait replaceOnce(\n  context,\n  \"generi\u0161e \u010dinjenicu C3 (Osoba je pod konstantnim stresom du\u017ei period)\",\n  \" i korisnik \u0107e se testirati za druge poreme\u0107aje iako stres nije visok\",\n  Word.InsertLocation.after\n);\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the content changes described by the diff:\n#  1. Typo fix: \"Percieved\" -> \"Perceived\" (Skala opazanog stresa / Perceived Stress Scale)\n#  2. Remove \" (Pravilo za 7)\" from the \"P6 (Pravilo za 7) - \" rule label.\n#  3. In rule P7:\n#       \"posljednjih 6 mjeseci\" -> \"posljednja 3 mjeseca\"\n#       \"osim niskog do umjerenog stresa, kao ni predispozicije za oboljenjima\"\n#           -> \"osim stresa i/ili predispozicije ka oboljenjima\"\n#       append \" i korisnik ce se testirati za druge poremecaje iako stres nije visok\"\n#       after the C3 sentence.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1. Fix the \"Percieved\" -> \"Perceived\" typo.\nReplace-Text \"Percieved\" \"Perceived\"\n\n# 2. Drop the \" (Pravilo za 7)\" aside after \"P6\".\nReplace-Text \" (Pravilo za 7)\" \"\"\n\n# 3a. Tighten the recency window in rule P7.\nReplace-Text \"posljednjih 6 mjeseci\" \"posljednja 3 mjeseca\"\n\n# 3b. Reword the exclusion clause in rule P7.\nReplace-Text \"osim niskog do umjerenog stresa, kao ni predispozicije za oboljenjima\" \"osim stresa i/ili predispozicije ka oboljenjima\"\n\n# 3c. Append the new trailing clause onto the end of the C3 sentence.\nReplace-Text \"generi\u0161e \u010dinjenicu C3 (Osoba je pod konstantnim stresom du\u017ei period)\" \"generi\u0161e \u010dinjenicu C3 (Osoba je pod konstantnim stresom du\u017ei period) i korisnik \u0107e se testirati za druge poreme\u0107aje iako stres nije visok\"\n"}
